# Applies the "cryptos list" price/volume/coin-order refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" values are stored as plain text (e.g. "61.686.57", "0.997") in
# this workbook, even when they look numeric. Excel normally auto-converts a plain
# numeric-looking string typed into a General-formatted cell into a real number, so
# for every Price cell we briefly mark it as Text, assign the literal string, then
# restore its original style -- this keeps the value as text without leaving any
# lasting formatting change on the cell.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '61.686.57'

# Row 3
Set-TextValue $ws.Range("D3") '3.399.14'
$ws.Range("E3").Value = '  -0.61%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.997'
$ws.Range("E4").Value = '  -0.26%  '

# Row 5
Set-TextValue $ws.Range("D5") '408.28'
$ws.Range("E5").Value = '  -0.32%  '

# Row 6
Set-TextValue $ws.Range("D6") '127.03'
$ws.Range("E6").Value = '  -1.45%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.613'
$ws.Range("E7").Value = '  -2.34%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.715'
$ws.Range("E9").Value = '  -4.55%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.133'
$ws.Range("E10").Value = '  -9.03%  '

# Row 11
Set-TextValue $ws.Range("D11") '42.08'
$ws.Range("E11").Value = '  -1.01%  '

# Row 12
$ws.Range("E12").Value = '  -0.36%  '

# Row 13
Set-TextValue $ws.Range("D13") '3.923.77'
$ws.Range("E13").Value = '  -0.91%  '

# Row 14
Set-TextValue $ws.Range("D14") '9.02'
$ws.Range("E14").Value = '  +0.62%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D15") '0.0000206'
$ws.Range("E15").Value = '  -7.32%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D16") '20.32'
$ws.Range("E16").Value = '  -3.46%  '

# Row 17
Set-TextValue $ws.Range("D17") '3.395.73'
$ws.Range("E17").Value = '  -0.94%  '

# Row 18
Set-TextValue $ws.Range("D18") '12.12'
$ws.Range("E18").Value = '  -2.20%  '

# Row 19
Set-TextValue $ws.Range("D19") '1.06'
$ws.Range("E19").Value = '  -0.30%  '

# Row 20
Set-TextValue $ws.Range("D20") '61.601.52'
$ws.Range("E20").Value = '  -0.61%  '

# Row 21
Set-TextValue $ws.Range("D21") '481.28'
$ws.Range("E21").Value = '  +19.99%  '

# Row 22
Set-TextValue $ws.Range("D22") '89.07'
$ws.Range("E22").Value = '  -0.60%  '

# Row 23
Set-TextValue $ws.Range("D23") '3.19'
$ws.Range("E23").Value = '  +0.35%  '

# Row 24
Set-TextValue $ws.Range("D24") '13.03'
$ws.Range("E24").Value = '  -1.02%  '

# Row 25
Set-TextValue $ws.Range("D25") '3.25'
$ws.Range("E25").Value = '  +0.79%  '

# Row 26
Set-TextValue $ws.Range("D26") '33.14'
$ws.Range("E26").Value = '  +1.40%  '

# Row 27
Set-TextValue $ws.Range("D27") '9.13'
$ws.Range("E27").Value = '  +5.10%  '

# Row 28
$ws.Range("E28").Value = '  +0.11%  '

# Row 29
Set-TextValue $ws.Range("D29") '7.87'
$ws.Range("E29").Value = '  +3.69%  '

# Row 30
Set-TextValue $ws.Range("D30") '2.73'
$ws.Range("E30").Value = '  +1.03%  '

# Row 31
Set-TextValue $ws.Range("D31") '11.72'
$ws.Range("E31").Value = '  -0.88%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.167'
$ws.Range("E32").Value = '  -3.21%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.111'
$ws.Range("E33").Value = '  -6.49%  '

# Row 34
Set-TextValue $ws.Range("D34") '40.66'
$ws.Range("E34").Value = '  -5.75%  '

# Row 35
$ws.Range("E35").Value = '  -0.75%  '

# Row 36
Set-TextValue $ws.Range("D36") '55.30'
$ws.Range("E36").Value = '  +2.93%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.0482'
$ws.Range("E37").Value = '  -3.18%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.999'
$ws.Range("E38").Value = '  +0.14%  '

# Row 39
$ws.Range("E39").Value = '  +4.14%  '

# Row 40
Set-TextValue $ws.Range("D40") '146.26'
$ws.Range("E40").Value = '  +3.29%  '

# Row 41
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D41") '0.133'
$ws.Range("E41").Value = '  -0.36%  '

# Row 42
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D42") '3.32'
$ws.Range("E42").Value = '  -1.37%  '

# Row 43
Set-TextValue $ws.Range("D43") '2.91'
$ws.Range("E43").Value = '  +0.20%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.05'
$ws.Range("E44").Value = '  +3.71%  '

# Row 45
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D45") '4.15'
$ws.Range("E45").Value = '  +0.77%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D46") '2.51'
$ws.Range("E46").Value = '  +4.14%  '

# Row 47
$ws.Range("E47").Value = '  +16.15%  '

# Row 48
Set-TextValue $ws.Range("D48") '16.26'
$ws.Range("E48").Value = '  -2.54%  '

# Row 49
Set-TextValue $ws.Range("D49") '0.143'
$ws.Range("E49").Value = '  +9.60%  '

# Row 50
Set-TextValue $ws.Range("D50") '21.70'
$ws.Range("E50").Value = '  +0.03%  '

# Row 51
Set-TextValue $ws.Range("D51") '112.17'
$ws.Range("E51").Value = '  +14.28%  '

